# Refresh the cryptos price table (Price / Volume(1h) columns, plus the
# Frax/EnergySwap row swap) to match the latest GitHub Actions scrape.
#
# Column D ("Price") is stored as text throughout the sheet (values such as
# "26.375.60" or leading-zero decimals are not valid numbers), so every
# Price write below is prefixed with a leading apostrophe. That is the
# standard Excel "force text" marker: it stops Excel from re-interpreting
# the literal as a number (which would silently drop significant trailing
# zeros, e.g. "4.950" -> 4.95, or flip tiny values into scientific
# notation) while keeping the cell's number format at "General", matching
# the original file.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.375.60"
$ws.Range("E2").Value = "  +0.97%  "
$ws.Range("D3").Value = "'1.688.36"
$ws.Range("E3").Value = "  +1.08%  "
$ws.Range("E4").Value = "  +0.84%  "
$ws.Range("D5").Value = "'218.66"
$ws.Range("E5").Value = "  +0.95%  "
$ws.Range("D6").Value = "'0.5454"
$ws.Range("E6").Value = "  +4.53%  "
$ws.Range("E7").Value = "  +0.93%  "
$ws.Range("D8").Value = "'0.2722"
$ws.Range("E8").Value = "  +0.98%  "
$ws.Range("D9").Value = "'0.06449"
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("D10").Value = "'22.02"
$ws.Range("E10").Value = "  +0.79%  "
$ws.Range("D11").Value = "'0.07695"
$ws.Range("E11").Value = "  +3.34%  "
$ws.Range("D12").Value = "'1.690.17"
$ws.Range("E12").Value = "  +1.08%  "
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("D14").Value = "'0.5806"
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("D15").Value = "'0.000008319"
$ws.Range("E15").Value = "  -2.35%  "
$ws.Range("D16").Value = "'65.15"
$ws.Range("E16").Value = "  +1.53%  "
$ws.Range("D17").Value = "'26.418.55"
$ws.Range("E17").Value = "  +1.97%  "
$ws.Range("D18").Value = "'4.950"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("E19").Value = "  +0.84%  "
$ws.Range("D20").Value = "'10.97"
$ws.Range("E20").Value = "  +1.52%  "
$ws.Range("D21").Value = "'190.30"
$ws.Range("E21").Value = "  +0.47%  "
$ws.Range("D22").Value = "'6.226"
$ws.Range("E22").Value = "  +0.46%  "
$ws.Range("D23").Value = "'1.012"
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("D24").Value = "'149.57"
$ws.Range("E24").Value = "  +3.37%  "
$ws.Range("D25").Value = "'0.1305"
$ws.Range("E25").Value = "  +4.87%  "
$ws.Range("D26").Value = "'7.879"
$ws.Range("E26").Value = "  +3.56%  "
$ws.Range("D27").Value = "'15.72"
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("D28").Value = "'0.06354"
$ws.Range("E28").Value = "  -3.32%  "
$ws.Range("D29").Value = "'1.416"
$ws.Range("E29").Value = "  +5.44%  "
$ws.Range("E30").Value = "  +0.95%  "
$ws.Range("D31").Value = "'3.574"
$ws.Range("E31").Value = "  -0.38%  "
$ws.Range("E32").Value = "  +1.18%  "
$ws.Range("D33").Value = "'1.672"
$ws.Range("E33").Value = "  +0.36%  "
$ws.Range("E34").Value = "  +2.16%  "
$ws.Range("D35").Value = "'0.6194"
$ws.Range("E35").Value = "  +0.34%  "
$ws.Range("D37").Value = "'2.718"
$ws.Range("E37").Value = "  +0.79%  "
$ws.Range("D38").Value = "'6.228"
$ws.Range("E38").Value = "  -0.45%  "
$ws.Range("D39").Value = "'1.115.67"
$ws.Range("E39").Value = "  +1.91%  "
$ws.Range("D40").Value = "'0.01634"
$ws.Range("D41").Value = "'0.8803"
$ws.Range("E41").Value = "  +1.27%  "
$ws.Range("E42").Value = "  +0.67%  "
$ws.Range("D43").Value = "'101.01"
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("D44").Value = "'1.841.24"
$ws.Range("E44").Value = "  +1.30%  "
$ws.Range("E45").Value = "  -1.20%  "
$ws.Range("D46").Value = "'57.24"
$ws.Range("E46").Value = "  +1.31%  "
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").Value = "'1.011"
$ws.Range("E47").Value = "  +0.65%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'8.190"
$ws.Range("E48").Value = "  +0.76%  "
$ws.Range("D49").Value = "'0.05275"
$ws.Range("E49").Value = "  +0.72%  "
$ws.Range("D50").Value = "'0.4308"
$ws.Range("E50").Value = "  +0.70%  "
$ws.Range("D51").Value = "'6.030"
$ws.Range("E51").Value = "  +0.56%  "
